$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new values look numeric,
# so Excel stores them as text (preserving exact digits/trailing zeros)
# instead of silently converting them to floating-point numbers.
$ws.Range("D2").Value = '27.664.20'
$ws.Range("E2").Value = '  +0.37%  '

$ws.Range("D3").Value = '1.846.74'
$ws.Range("E3").Value = '  +0.03%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.032'
$ws.Range("E4").Value = '  +0.26%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '320.16'
$ws.Range("E5").Value = '  +0.19%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.027'
$ws.Range("E6").Value = '  +0.03%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4379'
$ws.Range("E7").Value = '  +0.21%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3784'
$ws.Range("E8").Value = '  +1.05%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07383'
$ws.Range("E9").Value = '  -0.13%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8794'
$ws.Range("E10").Value = '  +0.37%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.50'

$ws.Range("D12").Value = '1.869.51'
$ws.Range("E12").Value = '  +1.13%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.493'

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.670'
$ws.Range("E14").Value = '  -0.28%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07148'
$ws.Range("E15").Value = '  +0.01%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '84.74'
$ws.Range("E16").Value = '  +2.48%  '

$ws.Range("E17").Value = '  +0.14%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009057'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.028'
$ws.Range("E19").Value = '  +0.16%  '

$ws.Range("E20").Value = '  +0.16%  '

$ws.Range("D21").Value = '27.679.07'
$ws.Range("E21").Value = '  +0.38%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.268'
$ws.Range("E22").Value = '  +0.07%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.23'
$ws.Range("E23").Value = '  +0.22%  '

$ws.Range("D24").Value = '2.096.39'
$ws.Range("E24").Value = '  +1.48%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.045'
$ws.Range("E25").Value = '  +5.96%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '158.43'
$ws.Range("E26").Value = '  +0.58%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.65'
$ws.Range("E27").Value = '  -0.37%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.986'
$ws.Range("E28").Value = '  +2.26%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.319'
$ws.Range("E29").Value = '  +1.15%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '117.56'
$ws.Range("E30").Value = '  +1.30%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09052'
$ws.Range("E31").Value = '  -0.34%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7724'
$ws.Range("E32").Value = '  +0.54%  '

$ws.Range("E33").Value = '  +0.17%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.991'
$ws.Range("E34").Value = '  +3.89%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.546'
$ws.Range("E35").Value = '  +1.09%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.028'
$ws.Range("E36").Value = '  +0.02%  '

$ws.Range("E37").Value = '  +0.32%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01970'
$ws.Range("E38").Value = '  -0.28%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05253'
$ws.Range("E39").Value = '  -0.14%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.841'
$ws.Range("E40").Value = '  +1.50%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5169'
$ws.Range("E41").Value = '  -0.11%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1666'
$ws.Range("E42").Value = '  -0.42%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.828'
$ws.Range("E43").Value = '  +2.28%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.697'
$ws.Range("E44").Value = '  +2.17%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '109.90'
$ws.Range("E45").Value = '  +0.95%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.73'
$ws.Range("E46").Value = '  +1.70%  '

$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.06569'
$ws.Range("E47").Value = '  +3.19%  '

$ws.Range("B48").Value = 'PaxDollar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.029'
$ws.Range("E48").Value = '  -0.12%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.704'
$ws.Range("E49").Value = '  -0.52%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4684'
$ws.Range("E50").Value = '  +0.72%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.881'
$ws.Range("E51").Value = '  -1.20%  '
